$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the title string in A1 (Jan 2008 -- October 2015) and drop the external workbook link
$ws.Cells.Item(1,1).Value = "Real Average Weekly Wages; Year/Year: January 2008 -- October 2015"

# 2. Strip the external-workbook formulas from column A (rows 6-98), keeping the same cached date values
$ws.Cells.Item(6,1).Value = 39448
$ws.Cells.Item(7,1).Value = 39479
$ws.Cells.Item(8,1).Value = 39508
$ws.Cells.Item(9,1).Value = 39539
$ws.Cells.Item(10,1).Value = 39569
$ws.Cells.Item(11,1).Value = 39600
$ws.Cells.Item(12,1).Value = 39630
$ws.Cells.Item(13,1).Value = 39661
$ws.Cells.Item(14,1).Value = 39692
$ws.Cells.Item(15,1).Value = 39722
$ws.Cells.Item(16,1).Value = 39753
$ws.Cells.Item(17,1).Value = 39783
$ws.Cells.Item(18,1).Value = 39814
$ws.Cells.Item(19,1).Value = 39845
$ws.Cells.Item(20,1).Value = 39873
$ws.Cells.Item(21,1).Value = 39904
$ws.Cells.Item(22,1).Value = 39934
$ws.Cells.Item(23,1).Value = 39965
$ws.Cells.Item(24,1).Value = 39995
$ws.Cells.Item(25,1).Value = 40026
$ws.Cells.Item(26,1).Value = 40057
$ws.Cells.Item(27,1).Value = 40087
$ws.Cells.Item(28,1).Value = 40118
$ws.Cells.Item(29,1).Value = 40148
$ws.Cells.Item(30,1).Value = 40179
$ws.Cells.Item(31,1).Value = 40210
$ws.Cells.Item(32,1).Value = 40238
$ws.Cells.Item(33,1).Value = 40269
$ws.Cells.Item(34,1).Value = 40299
$ws.Cells.Item(35,1).Value = 40330
$ws.Cells.Item(36,1).Value = 40360
$ws.Cells.Item(37,1).Value = 40391
$ws.Cells.Item(38,1).Value = 40422
$ws.Cells.Item(39,1).Value = 40452
$ws.Cells.Item(40,1).Value = 40483
$ws.Cells.Item(41,1).Value = 40513
$ws.Cells.Item(42,1).Value = 40544
$ws.Cells.Item(43,1).Value = 40575
$ws.Cells.Item(44,1).Value = 40603
$ws.Cells.Item(45,1).Value = 40634
$ws.Cells.Item(46,1).Value = 40664
$ws.Cells.Item(47,1).Value = 40695
$ws.Cells.Item(48,1).Value = 40725
$ws.Cells.Item(49,1).Value = 40756
$ws.Cells.Item(50,1).Value = 40787
$ws.Cells.Item(51,1).Value = 40817
$ws.Cells.Item(52,1).Value = 40848
$ws.Cells.Item(53,1).Value = 40878
$ws.Cells.Item(54,1).Value = 40909
$ws.Cells.Item(55,1).Value = 40940
$ws.Cells.Item(56,1).Value = 40969
$ws.Cells.Item(57,1).Value = 41000
$ws.Cells.Item(58,1).Value = 41030
$ws.Cells.Item(59,1).Value = 41061
$ws.Cells.Item(60,1).Value = 41091
$ws.Cells.Item(61,1).Value = 41122
$ws.Cells.Item(62,1).Value = 41153
$ws.Cells.Item(63,1).Value = 41183
$ws.Cells.Item(64,1).Value = 41214
$ws.Cells.Item(65,1).Value = 41244
$ws.Cells.Item(66,1).Value = 41275
$ws.Cells.Item(67,1).Value = 41306
$ws.Cells.Item(68,1).Value = 41334
$ws.Cells.Item(69,1).Value = 41365
$ws.Cells.Item(70,1).Value = 41395
$ws.Cells.Item(71,1).Value = 41426
$ws.Cells.Item(72,1).Value = 41456
$ws.Cells.Item(73,1).Value = 41487
$ws.Cells.Item(74,1).Value = 41518
$ws.Cells.Item(75,1).Value = 41548
$ws.Cells.Item(76,1).Value = 41579
$ws.Cells.Item(77,1).Value = 41609
$ws.Cells.Item(78,1).Value = 41640
$ws.Cells.Item(79,1).Value = 41671
$ws.Cells.Item(80,1).Value = 41699
$ws.Cells.Item(81,1).Value = 41730
$ws.Cells.Item(82,1).Value = 41760
$ws.Cells.Item(83,1).Value = 41791
$ws.Cells.Item(84,1).Value = 41821
$ws.Cells.Item(85,1).Value = 41852
$ws.Cells.Item(86,1).Value = 41883
$ws.Cells.Item(87,1).Value = 41913
$ws.Cells.Item(88,1).Value = 41944
$ws.Cells.Item(89,1).Value = 41974
$ws.Cells.Item(90,1).Value = 42005
$ws.Cells.Item(91,1).Value = 42036
$ws.Cells.Item(92,1).Value = 42064
$ws.Cells.Item(93,1).Value = 42095
$ws.Cells.Item(94,1).Value = 42125
$ws.Cells.Item(95,1).Value = 42156
$ws.Cells.Item(96,1).Value = 42186
$ws.Cells.Item(97,1).Value = 42217
$ws.Cells.Item(98,1).Value = 42248

# 3. Update row 98 (C98:BA98) with revised figures
$ws.Cells.Item(98,3).Value = 1.1005017013327831
$ws.Cells.Item(98,4).Value = 2.1306278272832739
$ws.Cells.Item(98,5).Value = 3.4696053285147563
$ws.Cells.Item(98,6).Value = -1.4526918822004753
$ws.Cells.Item(98,7).Value = 1.894983351928399
$ws.Cells.Item(98,8).Value = 0.74700110289772503
$ws.Cells.Item(98,9).Value = 2.8513267917191158
$ws.Cells.Item(98,10).Value = 6.1646266236050664
$ws.Cells.Item(98,11).Value = -5.684904347437711
$ws.Cells.Item(98,12).Value = 1.4497611988441521
$ws.Cells.Item(98,13).Value = -0.031705414087965629
$ws.Cells.Item(98,14).Value = -3.4426260724629052
$ws.Cells.Item(98,15).Value = 1.673531816142465
$ws.Cells.Item(98,16).Value = 1.6667312468045905
$ws.Cells.Item(98,17).Value = 0.71629687102618744
$ws.Cells.Item(98,18).Value = 3.1267199055027248
$ws.Cells.Item(98,19).Value = 0.071290253060657491
$ws.Cells.Item(98,20).Value = 4.563695160672947
$ws.Cells.Item(98,21).Value = -1.9086364013446206
$ws.Cells.Item(98,22).Value = 2.7999722138947103
$ws.Cells.Item(98,23).Value = -0.004334800042173066
$ws.Cells.Item(98,24).Value = 2.9855764779201559
$ws.Cells.Item(98,25).Value = 1.8543778500428487
$ws.Cells.Item(98,26).Value = 0.41437712885608741
$ws.Cells.Item(98,27).Value = -2.6519795651680766
$ws.Cells.Item(98,28).Value = 0.081967882980174836
$ws.Cells.Item(98,29).Value = 2.5326894631134462
$ws.Cells.Item(98,30).Value = 5.0756568510671753
$ws.Cells.Item(98,31).Value = 6.4064024621575291
$ws.Cells.Item(98,32).Value = 2.6559036258003519
$ws.Cells.Item(98,33).Value = 3.2279932806311606
$ws.Cells.Item(98,34).Value = -0.59087405653857894
$ws.Cells.Item(98,35).Value = 1.5170276511489749
$ws.Cells.Item(98,36).Value = 1.1924621450056274
$ws.Cells.Item(98,37).Value = -0.17736380203700494
$ws.Cells.Item(98,38).Value = 2.5440748920622096
$ws.Cells.Item(98,39).Value = -1.2932586777004671
$ws.Cells.Item(98,40).Value = 2.7497568290862908
$ws.Cells.Item(98,41).Value = 2.9095014233026415
$ws.Cells.Item(98,42).Value = 0.70025673522697107
$ws.Cells.Item(98,43).Value = 1.7281385442674255
$ws.Cells.Item(98,44).Value = 3.5340516634853874
$ws.Cells.Item(98,45).Value = 0.90683564915752368
$ws.Cells.Item(98,46).Value = -0.38870022789273728
$ws.Cells.Item(98,47).Value = 1.0454474201444657
$ws.Cells.Item(98,48).Value = 2.3737022812335042
$ws.Cells.Item(98,49).Value = 4.6355056937479464
$ws.Cells.Item(98,50).Value = 4.30470046454833
$ws.Cells.Item(98,51).Value = -0.66096357107067227
$ws.Cells.Item(98,52).Value = -0.27867177946022642
$ws.Cells.Item(98,53).Value = -3.4194571934963092

# 4. Add new row 99 (A99:BA99)
$ws.Cells.Item(99,1).Value = 42278
$ws.Cells.Item(99,2).Value = 2.0109720293641753
$ws.Cells.Item(99,3).Value = 3.875389073312093
$ws.Cells.Item(99,4).Value = 1.6713264370977241
$ws.Cells.Item(99,5).Value = 3.4988375535301701
$ws.Cells.Item(99,6).Value = 0.43302394287164692
$ws.Cells.Item(99,7).Value = 2.0792631730655082
$ws.Cells.Item(99,8).Value = 1.2589464957269563
$ws.Cells.Item(99,9).Value = 3.5758594401086343
$ws.Cells.Item(99,10).Value = 6.3838400348970579
$ws.Cells.Item(99,11).Value = -8.3058752488200565
$ws.Cells.Item(99,12).Value = 2.6486635728379273
$ws.Cells.Item(99,13).Value = 1.6996638775173722
$ws.Cells.Item(99,14).Value = -1.4022928484078947
$ws.Cells.Item(99,15).Value = 1.7020244649269953
$ws.Cells.Item(99,16).Value = 1.6524800621482143
$ws.Cells.Item(99,17).Value = 1.0015035399525174
$ws.Cells.Item(99,18).Value = 4.1900086022102334
$ws.Cells.Item(99,19).Value = 0.70079406355546436
$ws.Cells.Item(99,20).Value = 4.4161672817386926
$ws.Cells.Item(99,21).Value = -0.65807158333573323
$ws.Cells.Item(99,22).Value = 3.1233671692587106
$ws.Cells.Item(99,23).Value = 1.559358402569391
$ws.Cells.Item(99,24).Value = 3.7401779270501736
$ws.Cells.Item(99,25).Value = 2.7034601909267413
$ws.Cells.Item(99,26).Value = 2.4437636772824591
$ws.Cells.Item(99,27).Value = -0.95244836323532156
$ws.Cells.Item(99,28).Value = -0.44515591527256249
$ws.Cells.Item(99,29).Value = 1.485635963187746
$ws.Cells.Item(99,30).Value = 6.0250952034602738
$ws.Cells.Item(99,31).Value = 5.8004642924436753
$ws.Cells.Item(99,32).Value = 4.202902598420188
$ws.Cells.Item(99,33).Value = 3.161060117917744
$ws.Cells.Item(99,34).Value = -1.8565824674211004
$ws.Cells.Item(99,35).Value = 2.077786971076335
$ws.Cells.Item(99,36).Value = 3.1326224220425574
$ws.Cells.Item(99,37).Value = 0.37679809409480747
$ws.Cells.Item(99,38).Value = 2.5351355087312619
$ws.Cells.Item(99,39).Value = -0.71203965732984431
$ws.Cells.Item(99,40).Value = 3.3543683199632057
$ws.Cells.Item(99,41).Value = 3.6929845617981316
$ws.Cells.Item(99,42).Value = -0.2734690259049502
$ws.Cells.Item(99,43).Value = 2.2594997596382607
$ws.Cells.Item(99,44).Value = 4.6245140590792895
$ws.Cells.Item(99,45).Value = 2.1989038850612483
$ws.Cells.Item(99,46).Value = -0.31218107179246302
$ws.Cells.Item(99,47).Value = 1.044125285500326
$ws.Cells.Item(99,48).Value = 3.0328898921944356
$ws.Cells.Item(99,49).Value = 3.4764430084489697
$ws.Cells.Item(99,50).Value = 5.6330882589008358
$ws.Cells.Item(99,51).Value = 0.4140175531080646
$ws.Cells.Item(99,52).Value = 0.99756798909318467
$ws.Cells.Item(99,53).Value = -2.8344035384745156

# 5. Update the frozen-pane view / selected cell
$ws.Application.ActiveWindow.ScrollRow = 63
$ws.Range("C92").Select()
